# [createScaleInfo] Separate items of a scale and imputation information
#
# Previously, the "Items_der_Skala" column (D) held either the items of a
# scale (e.g. skala1) OR the imputation/plausible-value variable names
# (e.g. pv_pooled, pvkat_pooled). This script splits that information into
# two distinct columns: "Items_der_Skala" keeps only scale items, and a new
# "Imputationen" column (E) holds the imputation information.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the imputation-information column
$ws.Range("E1").Value = "Imputationen"

# Row 2: skala1 - has items already in D2, no imputations -> E2 blank
$ws.Range("E2").Value = ""

# Row 3: pv_pooled - its value was stored in D3 (items column) but is really
# imputation info, so move it over to the new E3 cell and blank out D3
$ws.Range("E3").Value = $ws.Range("D3").Value2
$ws.Range("D3").Value = ""

# Row 4: pvkat_pooled - same move as row 3
$ws.Range("E4").Value = $ws.Range("D4").Value2
$ws.Range("D4").Value = ""
